$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.1
$summary.Range("B4").Value = 0.21
$summary.Range("B6").Value = 161
$summary.Range("B8").Value = 57
$summary.Range("B9").Value = 45.34

# ---------------------------------------------------------------------------
# Strategy Status sheet - momentum strategy row (row 11)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C11").Value = 99.16
$status.Range("D11").Value = 41
$status.Range("E11").Value = -0.84
$status.Range("F11").Value = -0.84
$status.Range("G11").Value = 26.83

# ---------------------------------------------------------------------------
# All Trades sheet - close trade #161 (row 162) and append trades 190 & 191
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G162").Value = 0.82
$allTrades.Range("H162").Value = "CLOSED"
$allTrades.Range("I162").Value = -1.2048
$allTrades.Range("J162").Value = -0.01
$allTrades.Range("K162").Value = 99.16
$allTrades.Range("L162").Value = "early_exit"
$allTrades.Range("M162").Value = 0.19

# New trade #190 - momentum - row 191
$allTrades.Range("A191").Value = 190
$allTrades.Range("B191").Value = "'2026-02-18"
$allTrades.Range("C191").Value = "00:41:05"
$allTrades.Range("D191").Value = "momentum"
$allTrades.Range("E191").Value = "DOWN"
$allTrades.Range("F191").Value = 0.83
$allTrades.Range("H191").Value = "OPEN"
$allTrades.Range("I191").Value = 0
$allTrades.Range("J191").Value = 0
$allTrades.Range("K191").Value = 99.16712996249174
$allTrades.Range("M191").Value = 0
$allTrades.Range("N191").Value = 0
$allTrades.Range("O191").Value = 0
$allTrades.Range("P191").Value = 0.9
$allTrades.Range("Q191").Value = "Downward momentum: -45.109% over 10 samples"

# New trade #191 - MarketMaking - row 192
$allTrades.Range("A192").Value = 191
$allTrades.Range("B192").Value = "'2026-02-18"
$allTrades.Range("C192").Value = "00:41:05"
$allTrades.Range("D192").Value = "MarketMaking"
$allTrades.Range("E192").Value = "UP"
$allTrades.Range("F192").Value = 0.17
$allTrades.Range("H192").Value = "OPEN"
$allTrades.Range("I192").Value = 0
$allTrades.Range("J192").Value = 0
$allTrades.Range("K192").Value = 99.28858346467945
$allTrades.Range("M192").Value = 0
$allTrades.Range("N192").Value = 0
$allTrades.Range("O192").Value = 0
$allTrades.Range("P192").Value = 0.6
$allTrades.Range("Q192").Value = "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------------
# momentum sheet - close trade #161 (row 42) and append trade 190 (row 51)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

$momentum.Range("G42").Value = 0.82
$momentum.Range("H42").Value = "CLOSED"
$momentum.Range("I42").Value = -1.2048
$momentum.Range("J42").Value = -0.01
$momentum.Range("K42").Value = 99.16
$momentum.Range("P42").Value = "early_exit"
$momentum.Range("Q42").Value = 0.19

$momentum.Range("A51").Value = 190
$momentum.Range("B51").Value = "'2026-02-18"
$momentum.Range("C51").Value = "00:41:05"
$momentum.Range("D51").Value = "momentum"
$momentum.Range("E51").Value = "DOWN"
$momentum.Range("F51").Value = 0.83
$momentum.Range("H51").Value = "OPEN"
$momentum.Range("I51").Value = 0
$momentum.Range("J51").Value = 0
$momentum.Range("K51").Value = 99.16712996249174
$momentum.Range("L51").Value = 0
$momentum.Range("M51").Value = 0
$momentum.Range("N51").Value = 0.9
$momentum.Range("O51").Value = "Downward momentum: -45.109% over 10 samples"
$momentum.Range("Q51").Value = 0

# ---------------------------------------------------------------------------
# MarketMaking sheet - append trade 191 (row 81)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

$marketMaking.Range("A81").Value = 191
$marketMaking.Range("B81").Value = "'2026-02-18"
$marketMaking.Range("C81").Value = "00:41:05"
$marketMaking.Range("D81").Value = "MarketMaking"
$marketMaking.Range("E81").Value = "UP"
$marketMaking.Range("F81").Value = 0.17
$marketMaking.Range("H81").Value = "OPEN"
$marketMaking.Range("I81").Value = 0
$marketMaking.Range("J81").Value = 0
$marketMaking.Range("K81").Value = 99.28858346467945
$marketMaking.Range("L81").Value = 0
$marketMaking.Range("M81").Value = 0
$marketMaking.Range("N81").Value = 0.6
$marketMaking.Range("O81").Value = "Normal spread capture: 198 bps"
$marketMaking.Range("Q81").Value = 0
